$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Exigences")

# Insert two new columns before column AC (29) - shifts old AC:AI to AE:AK
$ws.Range("AC1:AD1").EntireColumn.Insert()

# Set the new header labels (row 1) - order matters for shared-string indices:
# "Statut publication" must be interned before "Commentaire" (45, 46)
$ws.Range("AD1").Value = "Statut publication"
$ws.Range("AC1").Value = "Commentaire"

# The column insert doesn't relocate the existing hyperlink anchors, so
# recreate them over the (now shifted-right-by-2) cells that hold them.
$link1 = "https://saas-ans02.henix.com/squash/requirement-workspace/requirement/null/content"
$link2 = "https://saas-ans02.henix.com/squash/test-case-workspace/test-case/0/content"
$link3 = "https://saas-ans02.henix.com/squash/requirement-workspace/requirement/20255/content"

$origStyleAF2 = $ws.Range("AF2").Style
$origStyleAG2 = $ws.Range("AG2").Style
$origStyleAH2 = $ws.Range("AH2").Style

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("AF2"), $link1)
$ws.Hyperlinks.Add($ws.Range("AG2"), $link2)
$ws.Hyperlinks.Add($ws.Range("AH2"), $link3)

# Adding a hyperlink re-styles the cell with the built-in "Hyperlink" style;
# restore the original (non-hyperlink) formatting that the template uses.
$ws.Range("AF2").Style = $origStyleAF2
$ws.Range("AG2").Style = $origStyleAG2
$ws.Range("AH2").Style = $origStyleAH2
